$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "back" values in row 2 (B2:U2), leaving formatting/style intact.
$ws.Range("B2:U2").ClearContents()

# Update the active selection to match the final saved state.
$ws.Range("U22").Select()
